$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.10031222304565
$ws.Range("C2").Value = 11.93688335548763
$ws.Range("D2").Value = 5.131882537702682
$ws.Range("E2").Value = 9.893039654352632
$ws.Range("F2").Value = 30.69742812745088
$ws.Range("I2").Value = 22.91209905447371
$ws.Range("M2").Value = 16.21238624863642
$ws.Range("N2").Value = 17.73480112790801
$ws.Range("B3").Value = 14.49656486284222
$ws.Range("C3").Value = 11.28452866778954
$ws.Range("D3").Value = 5.158064285910764
$ws.Range("E3").Value = 9.806446140188722
$ws.Range("F3").Value = 30.32580373260803
$ws.Range("I3").Value = 22.8648791195312
$ws.Range("M3").Value = 15.9208498256995
$ws.Range("N3").Value = 17.80463104271804
$ws.Range("B4").Value = 14.11715935678368
$ws.Range("C4").Value = 10.86731478749961
$ws.Range("D4").Value = 5.174906245129845
$ws.Range("E4").Value = 9.75613700604222
$ws.Range("F4").Value = 30.10653554502559
$ws.Range("I4").Value = 22.84271331921623
$ws.Range("M4").Value = 15.74425087070027
$ws.Range("N4").Value = 17.84945308448488
$ws.Range("B5").Value = 13.96062383440404
$ws.Range("C5").Value = 10.69331013941731
$ws.Range("D5").Value = 5.181962297801123
$ws.Range("E5").Value = 9.736371651987991
$ws.Range("F5").Value = 30.01952353422079
$ws.Range("I5").Value = 22.83539627480172
$ws.Range("M5").Value = 15.67299155904075
$ws.Range("N5").Value = 17.8682088108309
$ws.Range("B6").Value = 13.9345232355915
$ws.Range("C6").Value = 10.66418272942906
$ws.Range("D6").Value = 5.183145599991933
$ws.Range("E6").Value = 9.733134546498452
$ws.Range("F6").Value = 30.00521950941275
$ws.Range("I6").Value = 22.83428483949073
$ws.Range("M6").Value = 15.66120460603623
$ws.Range("N6").Value = 17.87135283214098
$ws.Range("B7").Value = 14.11505568818165
$ws.Range("C7").Value = 10.86498394410652
$ws.Range("D7").Value = 5.175000624660099
$ws.Range("E7").Value = 9.755867442278921
$ws.Range("F7").Value = 30.10535246115355
$ws.Range("I7").Value = 22.84260769526417
$ws.Range("M7").Value = 15.74328684994622
$ws.Range("N7").Value = 17.84970404399677
$ws.Range("B8").Value = 14.89409574379574
$ws.Range("C8").Value = 11.71552157644016
$ws.Range("D8").Value = 5.14075104569992
$ws.Range("E8").Value = 9.862599363734876
$ws.Range("F8").Value = 30.56750493762175
$ws.Range("I8").Value = 22.89440047680656
$ws.Range("M8").Value = 16.11142606986684
$ws.Range("N8").Value = 17.75847515111609
$ws.Range("B9").Value = 16.34256762864584
$ws.Range("C9").Value = 13.24396669252635
$ws.Range("D9").Value = 5.079661462064379
$ws.Range("E9").Value = 10.09376866516848
$ws.Range("F9").Value = 31.53978443110501
$ws.Range("I9").Value = 23.05010523997669
$ws.Range("M9").Value = 16.84798749576563
$ws.Range("N9").Value = 17.59497098773446
$ws.Range("B10").Value = 17.34669929272456
$ws.Range("C10").Value = 14.27425869691709
$ws.Range("D10").Value = 5.03847308811205
$ws.Range("E10").Value = 10.27574289181344
$ws.Range("F10").Value = 32.28775897434536
$ws.Range("I10").Value = 23.19734239285945
$ws.Range("M10").Value = 17.39219442637144
$ws.Range("N10").Value = 17.48416269468308
$ws.Range("B11").Value = 17.7884680547886
$ws.Range("C11").Value = 14.72170801620011
$ws.Range("D11").Value = 5.020536217572426
$ws.Range("E11").Value = 10.36088360206447
$ws.Range("F11").Value = 32.63387970454252
$ws.Range("I11").Value = 23.2713806232491
$ws.Range("M11").Value = 17.63923393746356
$ws.Range("N11").Value = 17.43576244757022
$ws.Range("B12").Value = 17.95345810135599
$ws.Range("C12").Value = 14.88802111071717
$ws.Range("D12").Value = 5.013859008319397
$ws.Range("E12").Value = 10.39343879107918
$ws.Range("F12").Value = 32.76566636058787
$ws.Range("I12").Value = 23.30042263300961
$ws.Range("M12").Value = 17.73261094910254
$ws.Range("N12").Value = 17.4177221890837
$ws.Range("B13").Value = 17.91802876287656
$ws.Range("C13").Value = 14.85234264833051
$ws.Range("D13").Value = 5.015291944566409
$ws.Range("E13").Value = 10.38641387772297
$ws.Range("F13").Value = 32.73725378515888
$ws.Range("I13").Value = 23.29412337915862
$ws.Range("M13").Value = 17.71250965352267
$ws.Range("N13").Value = 17.42159469009615
$ws.Range("B14").Value = 17.80208867152746
$ws.Range("C14").Value = 14.73545371180733
$ws.Range("D14").Value = 5.0199845724034
$ws.Range("E14").Value = 10.3635557847338
$ws.Range("F14").Value = 32.64470819100375
$ws.Range("I14").Value = 23.27374984097595
$ws.Range("M14").Value = 17.64692009366246
$ws.Range("N14").Value = 17.43427249954299
$ws.Range("B15").Value = 17.73076893303019
$ws.Range("C15").Value = 14.6634468360975
$ws.Range("D15").Value = 5.022873933597623
$ws.Range("E15").Value = 10.34959471301088
$ws.Range("F15").Value = 32.58811118004621
$ws.Range("I15").Value = 23.26140106349886
$ws.Range("M15").Value = 17.60671947622994
$ws.Range("N15").Value = 17.44207549093069
$ws.Range("B16").Value = 17.31751418241539
$ws.Range("C16").Value = 14.24458295866185
$ws.Range("D16").Value = 5.039661406295473
$ws.Range("E16").Value = 10.2702241005562
$ws.Range("F16").Value = 32.26524647311814
$ws.Range("I16").Value = 23.19264507798409
$ws.Range("M16").Value = 17.37603096345294
$ws.Range("N16").Value = 17.48736608671977
$ws.Range("B17").Value = 17.0600447010813
$ws.Range("C17").Value = 13.98212794280221
$ws.Range("D17").Value = 5.050164894820395
$ws.Range("E17").Value = 10.22211954678253
$ws.Range("F17").Value = 32.06859225904221
$ws.Range("I17").Value = 23.15226719862921
$ws.Range("M17").Value = 17.23430727059757
$ws.Range("N17").Value = 17.51566390535935
$ws.Range("B18").Value = 16.91054869217496
$ws.Range("C18").Value = 13.82917521644661
$ws.Range("D18").Value = 5.056281510802981
$ws.Range("E18").Value = 10.19467473561012
$ws.Range("F18").Value = 31.95604133489206
$ws.Range("I18").Value = 23.1297080719152
$ws.Range("M18").Value = 17.15274784033749
$ws.Range("N18").Value = 17.53212901304451
$ws.Range("B19").Value = 16.85969492433978
$ws.Range("C19").Value = 13.77704778134614
$ws.Range("D19").Value = 5.058365420992773
$ws.Range("E19").Value = 10.1854215608035
$ws.Range("F19").Value = 31.91803340894736
$ws.Range("I19").Value = 23.1221844488123
$ws.Range("M19").Value = 17.12512871338249
$ws.Range("N19").Value = 17.53773629107287
$ws.Range("B20").Value = 17.08759938987905
$ws.Range("C20").Value = 14.01027383262413
$ws.Range("D20").Value = 5.049038988221379
$ws.Range("E20").Value = 10.22721739081199
$ws.Range("F20").Value = 32.0894694481349
$ws.Range("I20").Value = 23.15649671611156
$ws.Range("M20").Value = 17.24939923441662
$ws.Range("N20").Value = 17.51263200402585
$ws.Range("B21").Value = 17.83620647049097
$ws.Range("C21").Value = 14.76987217563939
$ws.Range("D21").Value = 5.018603109551636
$ws.Range("E21").Value = 10.37026143482953
$ws.Range("F21").Value = 32.67187257529871
$ws.Range("I21").Value = 23.27970684548155
$ws.Range("M21").Value = 17.66619073768116
$ws.Range("N21").Value = 17.4305409127657
$ws.Range("B22").Value = 18.31201763722908
$ws.Range("C22").Value = 15.24806781961366
$ws.Range("D22").Value = 4.999382413771744
$ws.Range("E22").Value = 10.46556867344637
$ws.Range("F22").Value = 33.05663579176552
$ws.Range("I22").Value = 23.3660858471495
$ws.Range("M22").Value = 17.93754694864522
$ws.Range("N22").Value = 17.37856721350872
$ws.Range("B23").Value = 18.05933963668498
$ws.Range("C23").Value = 14.99453551045121
$ws.Range("D23").Value = 5.009579442744726
$ws.Range("E23").Value = 10.41454333093988
$ws.Range("F23").Value = 32.85094370068472
$ws.Range("I23").Value = 23.31945180583713
$ws.Range("M23").Value = 17.79284498202968
$ws.Range("N23").Value = 17.40615329304031
$ws.Range("B24").Value = 17.07514649355936
$ws.Range("C24").Value = 13.99755549858197
$ws.Range("D24").Value = 5.049547767583442
$ws.Range("E24").Value = 10.22491199598527
$ws.Range("F24").Value = 32.08002927852892
$ws.Range("I24").Value = 23.15458251143718
$ws.Range("M24").Value = 17.24257640124612
$ws.Range("N24").Value = 17.51400211508582
$ws.Range("B25").Value = 15.96051161177794
$ws.Range("C25").Value = 12.84629915081955
$ws.Range("D25").Value = 5.09553856102568
$ws.Range("E25").Value = 10.02900877599988
$ws.Range("F25").Value = 31.27035968232315
$ws.Range("I25").Value = 23.00219624338636
$ws.Range("M25").Value = 16.64777455586015
$ws.Range("N25").Value = 17.6375620629165
